$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.267.94'
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("D3").Value = '2.279.87'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.21%  '
$ws.Range("E7").Value = '  +3.16%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.577'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.86%  '
$ws.Range("E13").Value = '  +3.08%  '
$ws.Range("D14").Value = '2.625.91'
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.883'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.72%  '
$ws.Range("D17").Value = '2.291.49'
$ws.Range("E17").Value = '  +4.83%  '
$ws.Range("D18").Value = '44.179.79'
$ws.Range("E18").Value = '  +3.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.21%  '
$ws.Range("E20").Value = '  +4.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("E25").Value = '  +5.39%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +19.16%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0893'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.63'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '161.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.93%  '
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.77%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.56'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.63%  '
$ws.Range("E41").Value = '  +5.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +35.58%  '
$ws.Range("E43").Value = '  +3.52%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").Value = '1.812.06'
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.210'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '77.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.08%  '
